$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/percentage cells and multi-dot "price" strings that Excel cannot
# parse as numbers -- safe to assign directly as .Value without any numeric
# auto-conversion happening.
$ws.Range('D2').Value = '21.476.33'
$ws.Range('E2').Value = '  -2.71%  '
$ws.Range('D3').Value = '1.529.82'
$ws.Range('E3').Value = '  -1.73%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('E6').Value = '  -1.19%  '
$ws.Range('E7').Value = '  -1.28%  '
$ws.Range('E8').Value = '  -1.59%  '
$ws.Range('E9').Value = '  -2.79%  '
$ws.Range('E10').Value = '  -2.31%  '
$ws.Range('E11').Value = '  -1.80%  '
$ws.Range('E12').Value = '  +0.17%  '
$ws.Range('E13').Value = '  +0.49%  '
$ws.Range('E14').Value = '  -4.10%  '
$ws.Range('E15').Value = '  -1.86%  '
$ws.Range('D16').Value = '1.531.71'
$ws.Range('E16').Value = '  -1.39%  '
$ws.Range('E17').Value = '  -4.22%  '
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('E19').Value = '  -0.83%  '
$ws.Range('E20').Value = '  +0.11%  '
$ws.Range('E21').Value = '  -3.45%  '
$ws.Range('E23').Value = '  -4.68%  '
$ws.Range('E24').Value = '  +0.66%  '
$ws.Range('D25').Value = '21.484.91'
$ws.Range('E25').Value = '  -2.72%  '
$ws.Range('E26').Value = '  -3.80%  '
$ws.Range('E27').Value = '  -0.40%  '
$ws.Range('E28').Value = '  -1.67%  '
$ws.Range('E29').Value = '  -1.62%  '
$ws.Range('D30').Value = '1.703.10'
$ws.Range('E31').Value = '  -2.34%  '
$ws.Range('E32').Value = '  +3.37%  '
$ws.Range('E33').Value = '  -5.87%  '
$ws.Range('E34').Value = '  -4.03%  '
$ws.Range('E35').Value = '  -6.74%  '
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('E37').Value = '  -8.21%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('E38').Value = '  -3.96%  '
$ws.Range('B39').Value = 'Aptos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('E39').Value = '  +3.85%  '
$ws.Range('E40').Value = '  -4.22%  '
$ws.Range('E41').Value = '  -2.23%  '
$ws.Range('E42').Value = '  -3.33%  '
$ws.Range('E43').Value = '  +0.11%  '
$ws.Range('E44').Value = '  -2.12%  '
$ws.Range('E45').Value = '  -0.29%  '
$ws.Range('E46').Value = '  -1.44%  '
$ws.Range('E47').Value = '  -1.32%  '
$ws.Range('E48').Value = '  -1.31%  '
$ws.Range('E49').Value = '  +1.66%  '
$ws.Range('E50').Value = '  -3.19%  '
$ws.Range('E51').Value = '  -2.48%  '

# "Price" cells whose new text would be parsed by Excel as a plain number
# (e.g. "288.23"). Force the cell to Text format before assigning so the
# literal string is preserved instead of being converted to a float, then
# reset the style back to Normal so no stray number-format style lingers on
# the cell (matches the original file, which carries no style index here).
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '288.23'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3887'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3176'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.55'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07137'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.066'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.711'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.12'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.517'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001083'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06603'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '83.06'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.087'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.39'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.79'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.369'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.355'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '148.38'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.30'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.818'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '116.40'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.018'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9515'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07988'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '8.486'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.146'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.491'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02187'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '11.21'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05872'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2011'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.174'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5725'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.08'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.707'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5534'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.888'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.158'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '115.25'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06665'
$ws.Range('D51').Style = 'Normal'
